# The commit removes the "Imię i nazwisko" (Name Surname) values that had
# been typed into column C (rows 4-26) of the "Arkusz1" worksheet.
# Clearing those cells causes the dependent formulas (TYPE/IF in columns
# D and E, and the SUM in G2) to recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Remove the typed-in names from C4:C26 (this also drops the now-unused
# shared strings, matching the sharedStrings.xml cleanup in the diff).
$ws.Range("C4:C26").ClearContents()

# Row 7 previously auto-grew (wrap text) to fit a long name; once the
# text is gone the row should shrink back to the default height.
$ws.Rows.Item(7).EntireRow.AutoFit()

# Reflect the author's final cursor position on the sheet.
$ws.Range("G12").Select() | Out-Null
